$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46033
$ws.Range("B2").Value = 96.90000000000001
$ws.Range("C2").Value = 85.69
$ws.Range("D2").Value = 81.56
$ws.Range("E2").Value = 82.62
$ws.Range("F2").Value = 81.56999999999999
$ws.Range("G2").Value = 82.23999999999999
$ws.Range("H2").Value = 81.59999999999999
$ws.Range("I2").Value = 84.26000000000001
$ws.Range("J2").Value = 90.63
$ws.Range("K2").Value = 83.92
$ws.Range("L2").Value = 69.53
$ws.Range("M2").Value = 54.02
$ws.Range("N2").Value = 50.47
$ws.Range("O2").Value = 53.05
$ws.Range("P2").Value = 55.1
$ws.Range("Q2").Value = 58.14
$ws.Range("R2").Value = 76.84
$ws.Range("S2").Value = 92.64
$ws.Range("T2").Value = 103.14
$ws.Range("U2").Value = 102.14
$ws.Range("V2").Value = 99.02
$ws.Range("W2").Value = 96.54000000000001
$ws.Range("X2").Value = 92.52
$ws.Range("Y2").Value = 80.97
$ws.Range("Z2").Value = 80.63
$ws.Range("AB2").Value = 93.69
$ws.Range("AD2").Value = 102.64
$ws.Range("AF2").Value = 97.78
$ws.Range("AG2").Value = "10h-16h"
